$d = $word.ActiveDocument
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Titel"/></w:pPr><w:r><w:t xml:space="preserve">OC0101 - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>angivBruttofortjenesteOgVareforbrug</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/></w:pBdr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Overskrift1"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Operation:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>angivBruttofortjenesteOg</w:t></w:r><w:r><w:t>Vareforbrug</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>(</w:t></w:r><w:r><w:t xml:space="preserve">bruttofortjeneste, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>vareforbrug</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Overskrift1"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>Cross references:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">UC01 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Beregn</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>omsætning</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Overskrift1"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>Preconditions:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>En</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>HØK h</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>eksisterer</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">h </w:t></w:r><w:r><w:t>h</w:t></w:r><w:r><w:t>ar ingen associationer</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">bruttofortjeneste </w:t></w:r><w:r><w:t>er</w:t></w:r><w:r><w:t xml:space="preserve"> større end 0</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>vareforbrug</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>er større end 0</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Overskrift1"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Postconditions</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">omsætning </w:t></w:r><w:r><w:t>blev</w:t></w:r><w:r><w:t xml:space="preserve"> udregnet</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">omsætning </w:t></w:r><w:r><w:t>blev præsenteret for h</w:t></w:r></w:p>
'@
$d.Content.InsertXML($xml)
Write-Output $d.Paragraphs.Count
